# Updates cryptos list price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.538.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "264.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5243"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3233"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06811"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7810"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07777"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.843.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.027"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007989"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.576.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.639"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.478"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.031"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.188"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.683"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "111.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.188"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08729"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04838"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.132"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7214"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.876"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.117"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.08%  "

# Rows 37/38: VeChain and RenderToken swap places with refreshed data
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.260"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.84%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01791"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.4867"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9004"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "110.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.014"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.664"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4211"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.017"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1229"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8894"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.01%  "
